# Updated symbol list (crypto price/volume refresh) - applies the new
# Price (column D) and Volume(1h) (column E) values for the affected rows.
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the workbook's existing inlineStr/text representation) instead
# of coercing numeric- or percent-looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "318.11";   "E2"  = "4.84%"
    "D3"  = "48.61";    "E3"  = "13.27%"
    "D4"  = "5.272";    "E4"  = "4.64%"
    "D5"  = "0.08009";  "E5"  = "4.22%"
    "D6"  = "4.588";    "E6"  = "4.02%"
    "D7"  = "1.431";    "E7"  = "35.30%"
    "D8"  = "1.649";    "E8"  = "2.39%"
    "D9"  = "0.1278";   "E9"  = "3.76%"
    "D10" = "0.1945";   "E10" = "5.08%"
    "D11" = "0.09321";  "E11" = "3.06%"
    "D12" = "0.04593";  "E12" = "10.34%"
    "E13" = "0.16%"
    "D14" = "0.001327"; "E14" = "4.63%"
    "D15" = "0.04175";  "E15" = "0.58%"
    "D16" = "0.005859"; "E16" = "1.46%"
    "D17" = "3.334";    "E17" = "0.51%"
    "D18" = "2.441";    "E18" = "2.38%"
    "D19" = "0.3420";   "E19" = "2.26%"
    "D20" = "8.170"
    "E21" = "-0.26%"
    "D22" = "0.3098";   "E22" = "7.17%"
    "D23" = "0.001309"; "E23" = "2.78%"
    "D24" = "0.004244"; "E24" = "-5.38%"
    "D25" = "0.0001351"; "E25" = "0.33%"
    "D26" = "0.0003540"; "E26" = "-95.24%"
    "D38" = "0.02690";  "E38" = "9.73%"
    "D39" = "0.05697";  "E39" = "7.97%"
    "E40" = "6.50%"
    "D41" = "0.008004"; "E41" = "4.57%"
    "D42" = "0.1440";   "E42" = "6.96%"
    "D43" = "0.007683"; "E43" = "4.34%"
    "D44" = "0.007893"; "E44" = "-5.54%"
    "D45" = "0.3492";   "E45" = "14.31%"
    "D46" = "0.00006892"; "E46" = "4.18%"
    "E47" = "0.36%"
    "D48" = "0.05486";  "E48" = "43.22%"
    "D49" = "0.004001"; "E49" = "-4.75%"
    "E50" = "0.36%"
    "D51" = "0.0002000"; "E51" = "0.36%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
